$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H132").Value = 1718.2094
$ws.Range("I132").Value = 1672.7805
$ws.Range("K132").Value = 5018.3415
$ws.Range("M132").Value = -2488.3415
$ws.Range("H137").Value = 1693.0189
$ws.Range("I137").Value = 1387.7073
$ws.Range("J137").Value = 2736.1667
$ws.Range("K137").Value = 4163.1219
$ws.Range("L137").Value = 8208.500100000001
$ws.Range("M137").Value = -1613.1219
$ws.Range("N137").Value = -13308.5001
$ws.Range("H138").Value = 6472.9673
$ws.Range("I138").Value = 3579
$ws.Range("K138").Value = 10737
$ws.Range("M138").Value = -5597

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5162.8677
$ws.Range("J32").Value = 9999.5
$ws.Range("L32").Value = 9999.5
$ws.Range("N32").Value = -10573.5
$ws.Range("H45").Value = 2413.9333
$ws.Range("I45").Value = 2050.875
$ws.Range("J45").Value = 2828.8572
$ws.Range("K45").Value = 2050.875
$ws.Range("L45").Value = 2828.8572
$ws.Range("M45").Value = -1673.875
$ws.Range("N45").Value = -3582.8572
$ws.Range("H61").Value = 2122.5518
$ws.Range("I61").Value = 2487
$ws.Range("K61").Value = 2487
$ws.Range("M61").Value = -2275
$ws.Range("H74").Value = 1507.091
$ws.Range("I74").Value = 1650.9697
$ws.Range("J74").Value = 1075.4546
$ws.Range("K74").Value = 1650.9697
$ws.Range("L74").Value = 1075.4546
$ws.Range("M74").Value = -776.9697000000001
$ws.Range("N74").Value = -2823.4546
$ws.Range("H77").Value = 1507.091
$ws.Range("I77").Value = 1650.9697
$ws.Range("J77").Value = 1075.4546
$ws.Range("K77").Value = 8254.8485
$ws.Range("L77").Value = 5377.273
$ws.Range("M77").Value = -3886.8485
$ws.Range("N77").Value = -14113.273
$ws.Range("H122").Value = 3591.8838
$ws.Range("I122").Value = 2526.76
$ws.Range("J122").Value = 5071.222
$ws.Range("K122").Value = 7580.280000000001
$ws.Range("L122").Value = 15213.666
$ws.Range("M122").Value = -5130.280000000001
$ws.Range("N122").Value = -20113.666
$ws.Range("H132").Value = 2821.9827
$ws.Range("I132").Value = 2665.1
$ws.Range("J132").Value = 3802.5
$ws.Range("K132").Value = 7995.299999999999
$ws.Range("L132").Value = 11407.5
$ws.Range("M132").Value = -5465.299999999999
$ws.Range("N132").Value = -16467.5
$ws.Range("H136").Value = 2122.5518
$ws.Range("I136").Value = 2487
$ws.Range("K136").Value = 7461
$ws.Range("M136").Value = -4911

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2118.4443
$ws.Range("I99").Value = 2038.8125
$ws.Range("K99").Value = 2038.8125
$ws.Range("M99").Value = -540.8125
$ws.Range("H107").Value = 418745.53
$ws.Range("I107").Value = 1930.5264
$ws.Range("J107").Value = 2002642.6
$ws.Range("K107").Value = 1930.5264
$ws.Range("L107").Value = 2002642.6
$ws.Range("M107").Value = -10.52639999999997
$ws.Range("N107").Value = -2006482.6
$ws.Range("H134").Value = 33183.35
$ws.Range("I134").Value = 4015.7036
$ws.Range("K134").Value = 12047.1108
$ws.Range("M134").Value = -9512.110799999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 103959.9
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 103959.9
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 103959.9
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -104549.9
$ws.Range("H34").Value = 103959.9
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 103959.9
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 103959.9
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -104363.9
$ws.Range("H122").Value = 5830.3
$ws.Range("I122").Value = 5600.5
$ws.Range("K122").Value = 16801.5
$ws.Range("M122").Value = -14351.5
$ws.Range("H132").Value = 1899.3125
$ws.Range("I132").Value = 1550.72
$ws.Range("J132").Value = 3144.2856
$ws.Range("K132").Value = 4652.16
$ws.Range("L132").Value = 9432.856800000001
$ws.Range("M132").Value = -2122.16
$ws.Range("N132").Value = -14492.8568
$ws.Range("H134").Value = 272658.75
$ws.Range("I134").Value = 2451.3547
$ws.Range("K134").Value = 7354.0641
$ws.Range("M134").Value = -4819.0641
$ws.Range("H135").Value = 72618.164
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 119916.664
$ws.Range("J37").Value = 119916.664
$ws.Range("L37").Value = 359749.992
$ws.Range("N37").Value = -359973.992
$ws.Range("H86").Value = 420
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -5072
$ws.Range("H89").Value = 420
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 8100
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -19956
$ws.Range("H92").Value = 1110.9286
$ws.Range("I92").Value = 769.8889
$ws.Range("J92").Value = 1724.8
$ws.Range("K92").Value = 2309.6667
$ws.Range("L92").Value = 5174.4
$ws.Range("M92").Value = -1061.6667
$ws.Range("N92").Value = -7670.4
$ws.Range("H107").Value = 88605.25
$ws.Range("J107").Value = 151048.58
$ws.Range("L107").Value = 453145.74
$ws.Range("N107").Value = -456985.74
$ws.Range("H131").Value = 2953.9473
$ws.Range("I131").Value = 2424
$ws.Range("J131").Value = 3143.2144
$ws.Range("K131").Value = 7272
$ws.Range("L131").Value = 9429.643199999999
$ws.Range("M131").Value = -2232
$ws.Range("N131").Value = -19509.6432

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8880.322
$ws.Range("I70").Value = 6521.647
$ws.Range("J70").Value = 11744.429
$ws.Range("K70").Value = 6521.647
$ws.Range("L70").Value = 11744.429
$ws.Range("M70").Value = -6251.647
$ws.Range("N70").Value = -12284.429
$ws.Range("H73").Value = 8880.322
$ws.Range("I73").Value = 6521.647
$ws.Range("J73").Value = 11744.429
$ws.Range("K73").Value = 6521.647
$ws.Range("L73").Value = 11744.429
$ws.Range("M73").Value = -5585.647
$ws.Range("N73").Value = -13616.429
$ws.Range("H97").Value = 746.6286
$ws.Range("J97").Value = 733.8
$ws.Range("L97").Value = 733.8
$ws.Range("N97").Value = -1725.8
$ws.Range("H122").Value = 3254.8333
$ws.Range("I122").Value = 3254.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9764.499899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7314.499899999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 22736.72
$ws.Range("I132").Value = 2905.3865
$ws.Range("K132").Value = 8716.1595
$ws.Range("M132").Value = -6186.1595

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5307.7085
$ws.Range("I7").Value = 5038.8
$ws.Range("K7").Value = 5038.8
$ws.Range("M7").Value = -4926.8
$ws.Range("H40").Value = 5573.3486
$ws.Range("I40").Value = 5040.931
$ws.Range("K40").Value = 5040.931
$ws.Range("M40").Value = -4904.931
$ws.Range("H61").Value = 5703.143
$ws.Range("I61").Value = 6722.6665
$ws.Range("J61").Value = 3868
$ws.Range("K61").Value = 6722.6665
$ws.Range("L61").Value = 3868
$ws.Range("M61").Value = -6520.6665
$ws.Range("N61").Value = -4272
$ws.Range("H113").Value = 5703.143
$ws.Range("I113").Value = 6722.6665
$ws.Range("J113").Value = 3868
$ws.Range("K113").Value = 6722.6665
$ws.Range("L113").Value = 3868
$ws.Range("M113").Value = -4552.6665
$ws.Range("N113").Value = -8208
$ws.Range("H122").Value = 4295.5
$ws.Range("J122").Value = 3200
$ws.Range("L122").Value = 9600
$ws.Range("N122").Value = -14500
$ws.Range("H126").Value = 5307.7085
$ws.Range("I126").Value = 5038.8
$ws.Range("K126").Value = 15116.4
$ws.Range("M126").Value = -12646.4
$ws.Range("H132").Value = 3745.5454
$ws.Range("I132").Value = 3314.7144
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 9944.143199999999
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -7414.143199999999
$ws.Range("N132").Value = -18558.5
$ws.Range("H139").Value = 40331.777
$ws.Range("J139").Value = 40331.777
$ws.Range("L139").Value = 40331.777
$ws.Range("N139").Value = -50611.777

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1728.9375
$ws.Range("I126").Value = 1702.1666
$ws.Range("J126").Value = 1809.25
$ws.Range("K126").Value = 5106.4998
$ws.Range("L126").Value = 5427.75
$ws.Range("M126").Value = -2636.4998
$ws.Range("N126").Value = -10367.75
$ws.Range("H136").Value = 56787
$ws.Range("I136").Value = 2842.2258
$ws.Range("J136").Value = 335501.66
$ws.Range("K136").Value = 8526.6774
$ws.Range("L136").Value = 1006504.98
$ws.Range("M136").Value = -5976.6774
$ws.Range("N136").Value = -1011604.98
